$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Next period (release date)" notes for the rows affected by the
# ONS SOC coding issue (Employment by occupation) and the ONS qualification
# framework recoding issue (Highest qualification level by age and gender).
$ws.Range("D3").Value = "TBC depending on ONS resolving the SOC coding issue"
$ws.Range("D7").Value = "TBC depending on ONS recoding the qualification framework."

# Move the active selection to D8 (as last edited by the author).
$ws.Range("D8").Select()
